# Add a new package entry row (row 11) to Sheet1, mirroring the existing
# "DTDemo" rows (8-10) but for a ContentPackage of version 1.0.0 uploaded
# on 2026-02-04.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a copy of row 10 so the new row picks up the same formatting
# and cell types (text) as its neighbours, including the F column which
# stores the upload date as text rather than a real date value.
$ws.Range("A10:F10").Copy($ws.Range("A11:F11"))

# Now overwrite with the values for the new entry. Column D already holds
# "1.0.0" and column F already holds "2026-02-04" (copied verbatim from row
# 10), so leave those two cells alone - re-assigning the F value as a
# string would make Excel auto-convert it into a real date serial, which
# is not what the source data looks like.
$ws.Range("A11").Value = "DTDemo"
$ws.Range("B11").Value = "DTDemo"
$ws.Range("C11").Value = "DTDemo"
$ws.Range("E11").Value = "ContentPackage"
